$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C25").Value = 1010
$ws.Range("D25").Value = 6043682
$ws.Range("E25").Value = 932.666975308642
$ws.Range("G25").Value = 7.675906183368864
$ws.Range("H25").Value = 26.57333870034762
